# Apply "test P7 with -10 percent" re-run results across all sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "general": scalar summary values
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 581.5593397042228
$ws.Range("B4").Value = 0.01799988746643066
$ws.Range("B6").Value = 33.08933970422287
$ws.Range("B7").Value = 5.336665625650533
$ws.Range("B8").Value = 5.336665625650533
$ws.Range("B9").Value = 487.1
$ws.Range("B10").Value = 61.37

# ---------------------------------------------------------------
# Sheet "x": column j (B) reassigned
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 4
$ws.Range("B5").Value = 13
$ws.Range("B6").Value = 3
$ws.Range("B9").Value = 10
$ws.Range("B10").Value = 12
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 9
$ws.Range("B13").Value = 11
$ws.Range("B14").Value = 7

# ---------------------------------------------------------------
# Sheet "U": column t (B) reassigned
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("U")
$ws.Range("B4").Value = 3
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 3

# ---------------------------------------------------------------
# Sheet "TBar": column TBar (B) reassigned
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value = 22.61192465059683
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 24.04101472405137
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 24.38986999490162
$ws.Range("B9").Value = 22.01159140980468
$ws.Range("B11").Value = 24.76592070603971
$ws.Range("B12").Value = 20
$ws.Range("B13").Value = 27.87444125446785
$ws.Range("B14").Value = 27.05494035044573
$ws.Range("B15").Value = 30

# ---------------------------------------------------------------
# Sheet "y": rows 2-4 edited, two new rows (5,6) appended
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("y")
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 13
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 13
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 3
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 1
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 1

# ---------------------------------------------------------------
# Sheet "Q": column Q (C) for rows 7-71 all reassigned
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 313.6
$ws.Range("C8").Value = 331.91
$ws.Range("C9").Value = 320.63
$ws.Range("C10").Value = 336.425
$ws.Range("C11").Value = 315.985
$ws.Range("C12").Value = 81.47500000000072
$ws.Range("C13").Value = 80.68000000000072
$ws.Range("C14").Value = 84.71500000000073
$ws.Range("C15").Value = 80.43500000000073
$ws.Range("C16").Value = 87.34500000000074
$ws.Range("C17").Value = 128.3950000000001
$ws.Range("C18").Value = 116.7850000000001
$ws.Range("C19").Value = 119.6300000000001
$ws.Range("C20").Value = 117.9250000000002
$ws.Range("C21").Value = 126.7800000000001
$ws.Range("C22").Value = 72.6299999999995
$ws.Range("C23").Value = 80.0549999999995
$ws.Range("C24").Value = 82.31999999999948
$ws.Range("C25").Value = 83.9549999999995
$ws.Range("C26").Value = 80.8149999999995
$ws.Range("C27").Value = 62.63000000000022
$ws.Range("C28").Value = 70.92000000000021
$ws.Range("C29").Value = 67.65500000000021
$ws.Range("C30").Value = 66.84500000000021
$ws.Range("C31").Value = 65.41000000000022
$ws.Range("C32").Value = 235.775
$ws.Range("C33").Value = 229.025
$ws.Range("C34").Value = 213.42
$ws.Range("C35").Value = 226.76
$ws.Range("C36").Value = 221.56
$ws.Range("C37").Value = 203.655
$ws.Range("C38").Value = 214.32
$ws.Range("C39").Value = 207.36
$ws.Range("C40").Value = 217.27
$ws.Range("C41").Value = 200.18
$ws.Range("C42").Value = 140.5549999999989
$ws.Range("C43").Value = 159.2149999999989
$ws.Range("C44").Value = 142.1399999999989
$ws.Range("C45").Value = 147.7249999999989
$ws.Range("C46").Value = 139.7449999999989
$ws.Range("C47").Value = 226.0399999999994
$ws.Range("C48").Value = 247.1799999999994
$ws.Range("C49").Value = 221.8549999999994
$ws.Range("C50").Value = 238.4549999999994
$ws.Range("C51").Value = 224.4749999999994
$ws.Range("C52").Value = 120.5799999999991
$ws.Range("C53").Value = 129.6
$ws.Range("C54").Value = 129.3849999999991
$ws.Range("C55").Value = 127.5
$ws.Range("C56").Value = 118.3249999999991
$ws.Range("C57").Value = 226.0399999999994
$ws.Range("C58").Value = 247.1799999999994
$ws.Range("C59").Value = 221.8549999999994
$ws.Range("C60").Value = 238.4549999999994
$ws.Range("C61").Value = 224.4749999999994
$ws.Range("C62").Value = 235.775
$ws.Range("C63").Value = 229.025
$ws.Range("C64").Value = 213.42
$ws.Range("C65").Value = 226.76
$ws.Range("C66").Value = 221.56
$ws.Range("C67").Value = 313.6
$ws.Range("C68").Value = 331.91
$ws.Range("C69").Value = 320.63
$ws.Range("C70").Value = 336.425
$ws.Range("C71").Value = 315.985

# ---------------------------------------------------------------
# Sheet "R": column R (C) reassigned
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R")
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C12").Value = 38.6
$ws.Range("C13").Value = 56.91
$ws.Range("C14").Value = 45.63
$ws.Range("C15").Value = 61.425
$ws.Range("C16").Value = 40.985

# ---------------------------------------------------------------
# Sheet "L": column L (C) reassigned
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("L")
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("C22").Value = 7.25
$ws.Range("C23").Value = 5.4
$ws.Range("C24").Value = 4.755
$ws.Range("C25").Value = 5.8
$ws.Range("C26").Value = 7.48
$ws.Range("C27").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("C31").Value = 0

# ---------------------------------------------------------------
# Sheet "rho": rows 2-4 edited, two new rows (5,6) appended
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1

# ---------------------------------------------------------------
# Sheet "alpha": rows 2-4 edited, two new rows (5,6) appended
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("alpha")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1
